$d = $word.ActiveDocument

# The paragraph currently reads "... began his training as a music producer
# and sound engineer at the Fred Plaut ..." split (among others) into two
# runs: " music producer and" and " sound engineer". We need the result to
# read "... began his training as a sound engineer and music producer at
# the Fred Plaut ...", with those same two runs becoming " " and
# "sound engineer and music producer" respectively (instead of being
# collapsed into a single run).

# Locate the unique combined phrase spanning the two runs we need to edit.
$whole = $d.Content
$whole.Find.Execute(" music producer and sound engineer", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$wholeStart = $whole.Start

$oldRun1Text = " music producer and"
$oldRun2Text = " sound engineer"
$newRun1Text = " "
$newRun2Text = "sound engineer and music producer"

$run1Start = $wholeStart
$run1End = $run1Start + $oldRun1Text.Length
$run2Start = $run1End
$run2End = $run2Start + $oldRun2Text.Length

$run1 = $d.Range($run1Start, $run1End)
$run2 = $d.Range($run2Start, $run2End)

if ($run1.Text -ne $oldRun1Text) { throw "Unexpected text in run1: [$($run1.Text)]" }
if ($run2.Text -ne $oldRun2Text) { throw "Unexpected text in run2: [$($run2.Text)]" }

# Temporarily apply a genuine (non-default) formatting marker so the text
# edits below do not get silently coalesced into neighbouring runs that
# happen to share identical (default) run formatting.
$run2.Font.Bold = $true
$run1.Font.Bold = $true

# Apply the text changes. run2 (the later range) is edited first so that
# run1's offsets remain valid while it, in turn, is edited.
$run2.Text = $newRun2Text
$run1.Text = $newRun1Text

# Remove the temporary Bold marker again: run1 first, then run2 using
# freshly recomputed offsets (run1 shrank, which shifts run2 to the left,
# and Range objects in this runtime do not auto-track such shifts).
$run1.Font.Bold = $false

$delta = $newRun1Text.Length - $oldRun1Text.Length
$run2FixedStart = $run2Start + $delta
$run2Fixed = $d.Range($run2FixedStart, $run2FixedStart + $newRun2Text.Length)
$run2Fixed.Font.Bold = $false
